$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-385 all hold the same serial date value 45178
# (2023-09-09). Update every one of them to 45179 (2023-09-10).
for ($row = 2; $row -le 385; $row++) {
    $ws.Cells.Item($row, 3).Value = 45179
}
